# Add 2022-Q4 data:
#   1. Insert a new summary row at the top of "总计" for the 2022-Q4 quarter.
#   2. Insert a brand-new "2022-Q4" worksheet (right after "总计") containing
#      the per-fund holdings for that quarter, copying formatting from the
#      existing "2022-Q3" sheet so styles/margins stay consistent.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" summary sheet: insert new row 2 = 2022-Q4 / 21 / 1.73
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Range("A2").EntireRow.Insert()
$summary.Range("A2:D2").ClearFormats()

# Re-apply the same cell style used by the sibling data rows (copy format
# from A3, which still carries the original "index column" style).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 21
$summary.Cells.Item(2, 4).Value = 1.73

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet with per-fund holdings
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2022-Q3")
$srcSheet.Copy($null, $summary)

$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The source sheet has 37 data rows; we only need 22 (header + 21 funds).
$q4.Range("A23:H37").EntireRow.Delete()
$q4.Range("A2:H22").ClearContents()

$data = @(
    @(0,  "004374", "华泰保兴吉年丰混合A",                 "5.90",  "94.77", "4.81", "0.2838", 7),
    @(1,  "011189", "建信智汇优选一年持有期混合（MOM）",   "17.54", "54.46", "1.56", "0.2736", 4),
    @(2,  "006642", "华泰保兴吉年利定期开放混合",           "6.73",  "93.27", "3.22", "0.2167", 9),
    @(3,  "011481", "广发瑞锦一年定开混合",                 "2.63",  "91.08", "7.29", "0.1917", 2),
    @(4,  "011336", "兴全汇吉一年持有期混合A",               "15.09", "39.83", "1.26", "0.1901", 9),
    @(5,  "005904", "华泰保兴成长优选混合A",                 "3.52",  "70.92", "4.35", "0.1531", 2),
    @(6,  "003857", "前海开源周期优选灵活配置混合A",         "2.13",  "89.59", "4.91", "0.1046", 8),
    @(7,  "006377", "广发趋势动力灵活配置混合",               "2.57",  "87.77", "3.12", "0.0802", 3),
    @(8,  "004375", "华泰保兴吉年丰混合C",                   "1.34",  "94.77", "4.81", "0.0645", 7),
    @(9,  "012132", "华泰保兴价值成长混合A",                 "0.80",  "81.95", "5.96", "0.0477", 1),
    @(10, "003858", "前海开源周期优选灵活配置混合C",         "0.72",  "89.59", "4.91", "0.0354", 8),
    @(11, "014999", "华泰保兴吉年盈混合A",                   "0.84",  "84.95", "4.02", "0.0338", 7),
    @(12, "350007", "天治趋势精选混合",                       "0.39",  "93.83", "4.22", "0.0165", 4),
    @(13, "011337", "兴全汇吉一年持有期混合C",               "0.80",  "39.83", "1.26", "0.0101", 9),
    @(14, "004931", "华润元大价值优选混合C",                 "0.16",  "74.11", "4.28", "0.0068", 7),
    @(15, "004930", "华润元大价值优选混合A",                 "0.13",  "74.11", "4.28", "0.0056", 7),
    @(16, "005905", "华泰保兴成长优选混合C",                 "0.11",  "70.92", "4.35", "0.0048", 2),
    @(17, "012177", "华泰保兴价值成长混合C",                 "0.08",  "81.95", "5.96", "0.0048", 1),
    @(18, "007257", "凯石沣混合A",                           "0.08",  "73.69", "4.01", "0.0032", 2),
    @(19, "007258", "凯石沣混合C",                           "0.04",  "73.69", "4.01", "0.0016", 2),
    @(20, "015000", "华泰保兴吉年盈混合C",                   "0.02",  "84.95", "4.02", "0.0008", 7)
)

$r = 2
foreach ($row in $data) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

Write-Host "2022-Q4 sheet added"
